$wb = $excel.ActiveWorkbook

# --- 1. Overview sheet: update status text for the 22f4966d... row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

# --- 2. zh-cn sheet: update status + add Error Detail (column L) for row 3 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handback transform failed"
$wsZh.Range("L3").Value = "Handback file name: cj352vco.qal is different with handoff file name: 22f4966d-b6ca-4330-91e7-a79f05217cce.c08626f4154cd03a3cf254215f6865d0857c3dd6.zh-cn."

# --- 3. de-de sheet: update status + add Error Detail (column L) for row 3 ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handback transform failed"
$wsDe.Range("L3").Value = "Handback file name: cj352vco.qal is different with handoff file name: 22f4966d-b6ca-4330-91e7-a79f05217cce.c08626f4154cd03a3cf254215f6865d0857c3dd6.de-de."
